$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.710788488388062
$ws.Range("B1").Value = 2.708538293838501
$ws.Range("C1").Value = 1.842766404151917
$ws.Range("D1").Value = 1.631384968757629
$ws.Range("E1").Value = 1.586158871650696
